{"js": "// Update the reported model-parameter numbers in the four\n// \"So the model parameters: slope ... and coefficients ...\" paragraphs.\n// Each paragraph has three separate text runs that need updating:\n//   intercept value, first coefficient value, second coefficient value.\n//\n// We locate the runs with Body.search() (exact, case-sensitive match on\n// the OLD text) and replace their content in place with insertText(...,\n// Word.InsertLocation.replace). Doing the three searches in order of\n// decreasing string length avoids any possibility of a shorter old value\n// accidentally being a substring of a not-yet-updated longer one.\n\nasync function replaceAll(context, oldText, newText) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n// 1) \"Fit the data: Using Scikit-Learn library\" section\n//    and 2) \"Fit the data: Using Custom Library OLS\" section\n//    (both sections currently show the identical numbers, so a single\n//    search/replace pass updates both occurrences).\nawait replaceAll(context, \"1.0055 and coefficients\", \"0.974 and coefficients\");\nawait replaceAll(context, \"2.9608, and\", \"2.9594, and\");\nawait replaceAll(context, \"2.0217\", \"2.0135\");\n\n// 3) \"Fit the data: Using Gradient Descent\" section\nawait replaceAll(context, \"1.0053 and coefficients\", \"0.9733 and coefficients\");\nawait replaceAll(context, \"2.9596, and\", \"2.9582, and\");\nawait replaceAll(context, \"2.0208\", \"2.0125\");\n\n// 4) \"Fit the data: Using Stochastic Gradient Descent\" section\nawait replaceAll(context, \"array([1.0119364]) and coefficients\", \"array([0.97379059]) and coefficients\");\nawait replaceAll(context, \"array([2.95776785]), and\", \"array([2.97930829]), and\");\nawait replaceAll(context, \"array([2.02670283])\", \"array([2.00812544])\");\n", "ps1": "# Update the reported model-parameter numbers in the four\n# \"So the model parameters: slope ... and coefficients ...\" paragraphs\n# (Scikit-learn, Custom OLS, Gradient Descent and Stochastic Gradient\n# Descent sections). Each value lives in its own run inside the\n# paragraph, so a plain whole-document Find/Replace on the exact old\n# number (scoped tightly enough to be unambiguous) updates every place\n# it appears.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Everywhere($oldText, $newText) {\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.MatchSoundsLike = $false\n    $find.MatchAllWordForms = $false\n    $find.Execute([ref]$oldText, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$newText, 2) | Out-Null\n}\n\n# 1) & 2) \"Fit the data: Using scikit-learn Library\" and\n#          \"Fit the data: Using Custom Library OLS\" sections share the\n#          same (identical) reported numbers before the edit.\nReplace-Everywhere \"1.0055 and coefficients\" \"0.974 and coefficients\"\nReplace-Everywhere \"2.9608, and\" \"2.9594, and\"\nReplace-Everywhere \"2.0217\" \"2.0135\"\n\n# 3) \"Fit the data: Using Gradient Descent\" section\nReplace-Everywhere \"1.0053 and coefficients\" \"0.9733 and coefficients\"\nReplace-Everywhere \"2.9596, and\" \"2.9582, and\"\nReplace-Everywhere \"2.0208\" \"2.0125\"\n\n# 4) \"Fit the data: Using Stochastic Gradient Descent\" section\nReplace-Everywhere \"array([1.0119364]) and coefficients\" \"array([0.97379059]) and coefficients\"\nReplace-Everywhere \"array([2.95776785]), and\" \"array([2.97930829]), and\"\nReplace-Everywhere \"array([2.02670283])\" \"array([2.00812544])\"\n"}
